# Generate Report for Handback
#
# Refreshes the "Latest Handback DateTime" for the handback report's first
# data row (1e0e29bd-dff7-4254-b748-5baad58266f7) on both locale sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Column K on both sheets is "Latest Handback DateTime"; row 2 is the
# 1e0e29bd-dff7-4254-b748-5baad58266f7 file.
$wsZhCn.Range("K2").Value = "2016-10-20 01:02:35"
$wsDeDe.Range("K2").Value = "2016-10-20 01:02:53"
